$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("E4").Value = "username"
$ws.Range("E4").Font.Underline = -4142
try {
  $ws.Range("E4").Font.Color = $null
  Write-Host "set null ok"
} catch {
  Write-Host "set null err: $_"
}
